# Improve image processing memory leaks
# Append one new trailing data row to each of the 4 worksheets, matching
# the existing layout/formatting of the row directly above it.

$wb = $excel.ActiveWorkbook

function Add-DataRow {
    param($ws, $row, $timeValue, $b, $c, $d, $e, $f, $g, $h, $i)

    $prevRow = $row - 1

    # Column A: timestamp, carries the same date/time number format as the
    # cell above it (style index 2 in the original workbook).
    $aCell = $ws.Cells.Item($row, 1)
    $aCell.Value = $timeValue
    $aCell.NumberFormat = $ws.Cells.Item($prevRow, 1).NumberFormat

    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
    $ws.Cells.Item($row, 9).Value = $i
}

# --- Sheet 1: ROW50-FE-LIFTER --- new row 24
$ws1 = $wb.Worksheets.Item("ROW50-FE-LIFTER")
$g1 = [double]"5.68631262647114e+23"
Add-DataRow $ws1 24 45736.13258740741 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c," "0x01,0x7e" "0xe" 400 $g1 382 14

# --- Sheet 2: ROW50-MID-LIFTER --- new row 26
# Column G on this sheet stores its (too-large-for-double-precision) ID
# value as literal text, matching every other row in the column. A leading
# apostrophe forces Excel to keep the digit string as text instead of
# silently re-parsing it back into a (rounded) number.
$ws2 = $wb.Worksheets.Item("ROW50-MID-LIFTER")
$g2 = "'568631262647113771663628"
Add-DataRow $ws2 26 45736.10865740741 "0x01,0x90 " "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," "0x01,0x82" "0x19" 400 $g2 386 25

# --- Sheet 3: ROW11-FE-LIFTER --- new row 24
$ws3 = $wb.Worksheets.Item("ROW11-FE-LIFTER")
$g3 = [double]"5.68631262647114e+23"
Add-DataRow $ws3 24 45736.15396599537 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c," "0x01,0x7e" "0x14" 400 $g3 382 20

# --- Sheet 4: ROW11-MID-LIFTER --- new row 24
$ws4 = $wb.Worksheets.Item("ROW11-MID-LIFTER")
$g4 = [double]"5.68631262647114e+23"
Add-DataRow $ws4 24 45736.3010040162 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c," "0x01,0x82" "0x19" 400 $g4 386 25

Write-Host "Appended new trailing rows to all four worksheets."
